# Applies the "glorot normal initialiser in GRU" update to the Results_Summary workbook.
# Target sheet: "Sheet4" (Tools/Results_Summary.xlsx -> 14 Days Ahead summary table)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")
$ws.Activate()

# --- New experiment block (rows 14-16): raw MAE / RMSE / R results for the
# new run, with the "Average" column (N) computed the same way as the rest
# of the 14-days-ahead table (AVERAGE across the 4 years in J:M). ---

# Row 14 - MAE
$ws.Range("J14").Value = 1.68553034493036
$ws.Range("K14").Value = 1.5020625826862199
$ws.Range("L14").Value = 1.85347392895776
$ws.Range("M14").Value = 3.53088879060993
$ws.Range("N14").Formula = "=AVERAGE(J14:M14)"

# Row 15 - RMSE
$ws.Range("J15").Value = 2.89145108940507
$ws.Range("K15").Value = 2.41428861776037
$ws.Range("L15").Value = 2.8916910345628999
$ws.Range("M15").Value = 7.0112605606291396
$ws.Range("N15").Formula = "=AVERAGE(J15:M15)"

# Row 16 - R
$ws.Range("J16").Value = 0.91113849598193597
$ws.Range("K16").Value = 0.94647925984624404
$ws.Range("L16").Value = 0.92030712483090504
$ws.Range("M16").Value = 0.96357185661028399
$ws.Range("N16").Formula = "=AVERAGE(J16:M16)"

# Row 16 closes the new block, so it gets the same thick bottom rule as every
# other last-row-of-a-block ("R" row) in this table.
$ws.Range("J16:N16").Borders.Item(9).LineStyle = 1
$ws.Range("J16:N16").Borders.Item(9).Weight = -4138
$ws.Rows.Item(16).RowHeight = 17

# --- Follow-up notes jotted in column L for the next iterations ---
$ws.Range("L20").Value = "ADD STDDEV"
$ws.Range("L18").Value = "3yrs weighted avg, patience around 20"
$ws.Range("L21").Value = "Add normal distribition centred to 0 for weights"

# Leave the cursor where the author last left it.
$ws.Range("L21").Select()
